$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.115.04"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").Value = "3.383.50"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'572.54"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "'136.28"
$ws.Range("E6").Value = "  +9.88%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.382.45"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").Value = "'0.477"
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("D10").Value = "'7.58"
$ws.Range("E10").Value = "  +4.44%  "
$ws.Range("D11").Value = "'0.123"
$ws.Range("E11").Value = "  +3.84%  "
$ws.Range("E12").Value = "  +4.05%  "
$ws.Range("D13").Value = "3.965.08"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").Value = "'0.0000174"
$ws.Range("E15").Value = "  +2.45%  "
$ws.Range("D16").Value = "3.392.42"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").Value = "'25.17"
$ws.Range("E17").Value = "  +3.70%  "
$ws.Range("D18").Value = "61.363.65"
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("D19").Value = "'14.03"
$ws.Range("E19").Value = "  +7.88%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'5.79"
$ws.Range("E20").Value = "  +3.77%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'9.42"
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("D22").Value = "'373.81"
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("D23").Value = "'0.568"
$ws.Range("E23").Value = "  +3.15%  "
$ws.Range("D24").Value = "3.519.36"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").Value = "'70.64"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("E27").Value = "  +12.65%  "
$ws.Range("D28").Value = "'1.67"
$ws.Range("E28").Value = "  +23.00%  "
$ws.Range("D29").Value = "'7.72"
$ws.Range("E29").Value = "  +13.21%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "'8.12"
$ws.Range("E31").Value = "  +5.66%  "
$ws.Range("E32").Value = "  +1.81%  "
$ws.Range("D33").Value = "'0.155"
$ws.Range("E33").Value = "  +4.80%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "3.415.92"
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").Value = "'23.43"
$ws.Range("E36").Value = "  +4.01%  "
$ws.Range("D37").Value = "'5.57"
$ws.Range("E37").Value = "  +9.49%  "
$ws.Range("D38").Value = "'1.57"
$ws.Range("E38").Value = "  +7.39%  "
$ws.Range("D39").Value = "'6.93"
$ws.Range("E39").Value = "  +5.42%  "
$ws.Range("D40").Value = "'163.09"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").Value = "'0.0789"
$ws.Range("E41").Value = "  +6.25%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'4.42"
$ws.Range("E43").Value = "  +5.12%  "
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").Value = "'1.21"
$ws.Range("E44").Value = "  +14.45%  "
$ws.Range("D45").Value = "'0.761"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").Value = "'41.38"
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("E47").Value = "  +5.70%  "
$ws.Range("D48").Value = "'23.24"
$ws.Range("E48").Value = "  +3.85%  "
$ws.Range("D49").Value = "'6.98"
$ws.Range("E49").Value = "  +6.14%  "
$ws.Range("D50").Value = "'23.00"
$ws.Range("E50").Value = "  +15.29%  "
$ws.Range("D51").Value = "'0.898"
$ws.Range("E51").Value = "  +7.14%  "
